$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 80: adenosine / nucleobase / oxygen / 1 / PMC7527729
$ws.Cells.Item(80, 1).Value = "adenosine"
$ws.Cells.Item(80, 2).Value = "nucleobase"
$ws.Cells.Item(80, 3).Value = "oxygen"
$ws.Cells.Item(80, 4).Value = 1
$ws.Cells.Item(80, 5).Value = "PMC7527729"

# Row 81: glucosamine / amino sugar / oxygen / 1 / PMC7527729
$ws.Cells.Item(81, 1).Value = "glucosamine"
$ws.Cells.Item(81, 2).Value = "amino sugar"
$ws.Cells.Item(81, 3).Value = "oxygen"
$ws.Cells.Item(81, 4).Value = 1
$ws.Cells.Item(81, 5).Value = "PMC7527729"

# Rows 82-85: the four deoxynucleotides, filled class-column-first (as a
# fill-down of "nucleotide" across B82:B85), then names in column A.
$ws.Cells.Item(82, 2).Value = "nucleotide"
$ws.Cells.Item(83, 2).Value = "nucleotide"
$ws.Cells.Item(84, 2).Value = "nucleotide"
$ws.Cells.Item(85, 2).Value = "nucleotide"

$ws.Cells.Item(82, 1).Value = "dAMP"
$ws.Cells.Item(83, 1).Value = "dTMP"
$ws.Cells.Item(84, 1).Value = "dGMP"
$ws.Cells.Item(85, 1).Value = "dCMP"

foreach ($r in 82..85) {
    $ws.Cells.Item($r, 3).Value = "oxygen"
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 35836416
}

# Rows 86-89: organic acids grown on oxygen.
$organicAcidRows = @(
    @(86, "pyruvate"),
    @(87, "fumarate"),
    @(88, "acetate"),
    @(89, "succinate")
)
foreach ($item in $organicAcidRows) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = "organic acid"
    $ws.Cells.Item($r, 3).Value = "oxygen"
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 16944129
}

# Rows 90-92: amino acids grown via fermentation.
$aminoAcidRows = @(
    @(90, "cysteine"),
    @(91, "serine"),
    @(92, "threonine")
)
foreach ($item in $aminoAcidRows) {
    $r = $item[0]
    $ws.Cells.Item($r, 1).Value = $item[1]
    $ws.Cells.Item($r, 2).Value = "amino acid"
    $ws.Cells.Item($r, 3).Value = "fermentation"
    $ws.Cells.Item($r, 4).Value = 1
    $ws.Cells.Item($r, 5).Value = 35377168
}

# Match the scrolled / selected view position recorded in the edited file.
$ws.Application.ActiveWindow.ScrollRow = 63
$ws.Range("E73").Select()
